$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All values in this sheet are plain text (inlineStr) cells in the source
# workbook, even when they look like numbers (e.g. "138.70", "0.0000231").
# Force text storage via NumberFormat "@" so Excel does not coerce these
# into numeric/scientific values and lose trailing zeros, then restore the
# default "Normal" style so no stray cell formatting is introduced.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.512.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.003.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.46%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "563.77"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  +0.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.988.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000231"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("E15").Value = "  +1.68%  "
$ws.Range("E16").Value = "  +6.51%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.497.91"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.996.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.406.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "431.25"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.721"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.40%  "
$ws.Range("E28").Value = "  +7.05%  "
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.16"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0990"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -7.84%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.994"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.60%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.92"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.72%  "
$ws.Range("B36").Value = "PEPE"
$ws.Range("C36").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0768"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.10"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "49.11"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.73"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "403.80"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.57%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0353"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.06%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.770.86"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.109"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.252"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.00%  "
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "120.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.63%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.57"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.44%  "
